$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.532.54"
$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "3.466.00"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.80%  "

$ws.Range("E7").Value = "  +1.79%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000277"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "4.019.55"
$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").Value = "3.460.43"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").Value = "65.489.01"
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "410.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.78%  "

$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "590.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.65%  "

$ws.Range("E33").Value = "  -1.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "60.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  -3.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.381"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("D41").Value = "3.214.29"
$ws.Range("E41").Value = "  +5.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.92%  "

$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("E45").Value = "  -5.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("E48").Value = "  -5.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.59%  "

$ws.Range("E51").Value = "  -2.67%  "
